# Update the Orthopedics United roster sheet with the latest
# player/position/team assignments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Chris Paul",         "PG",       "San Antonio Spurs"),
    @("Stephon Castle",     "PG,SG",    "San Antonio Spurs"),
    @("Payton Pritchard",   "PG,SG",    "Boston Celtics"),
    @("Brandin Podziemski", "SG",       "Golden State Warriors"),
    @("Deni Avdija",        "SF,PF",    "Portland Trail Blazers"),
    @("Paolo Banchero",     "SF,PF",    "Orlando Magic"),
    @("Chet Holmgren",      "PF,C",     "Oklahoma City Thunder"),
    @("Taurean Prince",     "SG,SF",    "Milwaukee Bucks"),
    @("Nikola Jokic",       "C",        "Denver Nuggets"),
    @("Pascal Siakam",      "SF,PF,C",  "Indiana Pacers"),
    @("Ayo Dosunmu",        "PG,SG,SF", "Chicago Bulls"),
    @("Rudy Gobert",        "C",        "Minnesota Timberwolves"),
    @("Jaylen Brown",       "SG,SF",    "Boston Celtics"),
    @("Jalen Green",        "PG,SG",    "Houston Rockets"),
    @("Dillon Brooks",      "SG,SF,PF", "Houston Rockets"),
    @("Russell Westbrook",  "PG,SG",    "Denver Nuggets"),
    @("Jalen Suggs",        "PG,SG",    "Orlando Magic"),
    @("Jakob Poeltl",       "C",        "Toronto Raptors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
